# Adds a second worksheet ("CalendarAsPrimary") with test-script data,
# mirroring the header formatting used on the "Report" sheet, and updates
# the selections on both sheets.

$wb = $excel.ActiveWorkbook
$report = $wb.Worksheets.Item("Report")

# --- Add the new sheet after "Report" ---
$calendar = $wb.Worksheets.Add($null, $report)
$calendar.Name = "CalendarAsPrimary"

# --- Header row ---
$calendar.Range("A1").Value = "ReportName"
$calendar.Range("B1").Value = "ReportFolder"
$calendar.Range("C1").Value = "PrimaryModule"
$calendar.Range("D1").Value = "ConditionField"
$calendar.Range("E1").Value = "CompareField"
$calendar.Range("F1").Value = "GroupField"

# --- Data row ---
$calendar.Range("A2").Value = "Chart_Report"
$calendar.Range("B2").Value = "Test Chethana"
$calendar.Range("C2").Value = "Calendar"
$calendar.Range("D2").Value = "Subject"
$calendar.Range("E2").Value = "Test"
$calendar.Range("F2").Value = "Subject"

# Re-use the same header styling (fill) as the "Report" sheet's header row.
[void]$report.Range("A1:F1").Copy()
[void]$calendar.Range("A1:F1").PasteSpecial(-4122)

# --- Selections ---
[void]$calendar.Range("F2").Select()
[void]$report.Range("A1:G2").Select()
[void]$report.Activate()
